$d = $word.ActiveDocument

function Replace-InRange($range, [string]$findText, [string]$replaceText) {
    # Locate-only Find (Replace=0 / wdReplaceNone) so the match isn't
    # rewritten by Find's own replace engine (which rebuilds/merges the
    # surrounding runs). Then duplicate the found (collapsed) range and
    # assign .Text on the duplicate, which only rewrites the matched
    # run's text node in place and leaves sibling runs untouched.
    $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $dup = $range.Duplicate
        $dup.Text = $replaceText
    }
    return $found
}

# --- Body text: "QWREW" -> "QWR" (bold run after "A ") ---
Replace-InRange $d.Content "QWREW" "QWR"

# --- Header replacements ---
# NOTE: re-fetch a fresh full-header Range before every call since a
# Range used by Find.Execute collapses down to the matched text.

# 1) standalone "QWREW" -> "QWR"  (must run BEFORE the "REW" replace,
#    since "QWREW" textually contains "REW")
Replace-InRange ($d.Sections.Item(1).Headers.Item(1).Range) "QWREW" "QWR"

# 2) "DIRETORIA DE ENSINO REGIAO REW" -> "...QWER"
Replace-InRange ($d.Sections.Item(1).Headers.Item(1).Range) "REW" "QWER"

# 3) "Rew" (5 occurrences in the address line) -> "Qwer"
for ($i = 0; $i -lt 5; $i++) {
    Replace-InRange ($d.Sections.Item(1).Headers.Item(1).Range) "Rew" "Qwer"
}

# 4-6) "rew" (3 occurrences: CEP, Tel, Email) -> "qwer"
for ($i = 0; $i -lt 3; $i++) {
    Replace-InRange ($d.Sections.Item(1).Headers.Item(1).Range) "rew" "qwer"
}
